$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1977.2933  # H15
$ws.Cells.Item(15, 9).Value = 1977.2933  # I15
$ws.Cells.Item(15, 11).Value = 5931.8799  # K15
$ws.Cells.Item(15, 13).Value = -5762.8799  # M15
$ws.Cells.Item(21, 8).Value = 4508.5  # H21
$ws.Cells.Item(21, 9).Value = 4508.5  # I21
$ws.Cells.Item(21, 11).Value = 4508.5  # K21
$ws.Cells.Item(21, 13).Value = -4040.5  # M21
$ws.Cells.Item(23, 8).Value = 4508.5  # H23
$ws.Cells.Item(23, 9).Value = 4508.5  # I23
$ws.Cells.Item(23, 11).Value = 4508.5  # K23
$ws.Cells.Item(23, 13).Value = -4274.5  # M23
$ws.Cells.Item(32, 8).Value = 950  # H32
$ws.Cells.Item(32, 9).Value = 950  # I32
$ws.Cells.Item(32, 10).Value = 0  # J32
$ws.Cells.Item(32, 11).Value = 950  # K32
$ws.Cells.Item(32, 12).Value = 0  # L32
$ws.Cells.Item(32, 13).Value = -624  # M32
$ws.Cells.Item(32, 14).ClearContents()  # N32
$ws.Cells.Item(33, 8).Value = 226  # H33
$ws.Cells.Item(33, 9).Value = 226  # I33
$ws.Cells.Item(33, 11).Value = 226  # K33
$ws.Cells.Item(33, 13).Value = 3  # M33
$ws.Cells.Item(40, 8).Value = 1306.3  # H40
$ws.Cells.Item(40, 9).Value = 750  # I40
$ws.Cells.Item(40, 10).Value = 1677.1666  # J40
$ws.Cells.Item(40, 11).Value = 750  # K40
$ws.Cells.Item(40, 12).Value = 1677.1666  # L40
$ws.Cells.Item(40, 13).Value = -575  # M40
$ws.Cells.Item(40, 14).Value = -2027.1666  # N40
$ws.Cells.Item(86, 8).Value = 50921  # H86
$ws.Cells.Item(86, 9).Value = 1800  # I86
$ws.Cells.Item(86, 10).Value = 100042  # J86
$ws.Cells.Item(86, 11).Value = 1800  # K86
$ws.Cells.Item(86, 12).Value = 100042  # L86
$ws.Cells.Item(86, 13).Value = -677  # M86
$ws.Cells.Item(86, 14).Value = -102288  # N86
$ws.Cells.Item(89, 8).Value = 50921  # H89
$ws.Cells.Item(89, 9).Value = 1800  # I89
$ws.Cells.Item(89, 10).Value = 100042  # J89
$ws.Cells.Item(89, 11).Value = 9000  # K89
$ws.Cells.Item(89, 12).Value = 500210  # L89
$ws.Cells.Item(89, 13).Value = -3384  # M89
$ws.Cells.Item(89, 14).Value = -511442  # N89
$ws.Cells.Item(132, 8).Value = 3498.5938  # H132
$ws.Cells.Item(132, 9).Value = 3459.6924  # I132
$ws.Cells.Item(132, 10).Value = 3667.1667  # J132
$ws.Cells.Item(132, 11).Value = 10379.0772  # K132
$ws.Cells.Item(132, 12).Value = 11001.5001  # L132
$ws.Cells.Item(132, 13).Value = -7849.0772  # M132
$ws.Cells.Item(132, 14).Value = -16061.5001  # N132
$ws.Cells.Item(137, 8).Value = 2284.6924  # H137
$ws.Cells.Item(137, 9).Value = 1974.8334  # I137
$ws.Cells.Item(137, 11).Value = 5924.5002  # K137
$ws.Cells.Item(137, 13).Value = -3374.5002  # M137
$ws.Cells.Item(138, 8).Value = 25003078  # H138
$ws.Cells.Item(138, 9).Value = 55557030  # I138
$ws.Cells.Item(138, 11).Value = 166671090  # K138
$ws.Cells.Item(138, 13).Value = -166665950  # M138

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1476.7222  # H2
$ws.Cells.Item(2, 9).Value = 1314.4375  # I2
$ws.Cells.Item(2, 11).Value = 1314.4375  # K2
$ws.Cells.Item(2, 13).Value = -1201.4375  # M2
$ws.Cells.Item(32, 8).Value = 6170.976  # H32
$ws.Cells.Item(32, 9).Value = 4893.5835  # I32
$ws.Cells.Item(32, 11).Value = 4893.5835  # K32
$ws.Cells.Item(32, 13).Value = -4606.5835  # M32
$ws.Cells.Item(36, 8).Value = 2513  # H36
$ws.Cells.Item(36, 9).Value = 2026  # I36
$ws.Cells.Item(36, 10).Value = 3000  # J36
$ws.Cells.Item(36, 11).Value = 2026  # K36
$ws.Cells.Item(36, 12).Value = 3000  # L36
$ws.Cells.Item(36, 13).Value = -1680  # M36
$ws.Cells.Item(36, 14).Value = -3692  # N36
$ws.Cells.Item(45, 8).Value = 3380.5715  # H45
$ws.Cells.Item(45, 9).Value = 3057.1428  # I45
$ws.Cells.Item(45, 11).Value = 3057.1428  # K45
$ws.Cells.Item(45, 13).Value = -2680.1428  # M45
$ws.Cells.Item(47, 8).Value = 24000  # H47
$ws.Cells.Item(47, 10).Value = 24000  # J47
$ws.Cells.Item(47, 12).Value = 24000  # L47
$ws.Cells.Item(47, 14).Value = -25450  # N47
$ws.Cells.Item(61, 8).Value = 3348.963  # H61
$ws.Cells.Item(61, 9).Value = 3216.88  # I61
$ws.Cells.Item(61, 10).Value = 5000  # J61
$ws.Cells.Item(61, 11).Value = 3216.88  # K61
$ws.Cells.Item(61, 12).Value = 5000  # L61
$ws.Cells.Item(61, 13).Value = -3004.88  # M61
$ws.Cells.Item(61, 14).Value = -5424  # N61
$ws.Cells.Item(97, 8).Value = 125001490  # H97
$ws.Cells.Item(97, 9).Value = 2690  # I97
$ws.Cells.Item(97, 10).Value = 250000290  # J97
$ws.Cells.Item(97, 11).Value = 2690  # K97
$ws.Cells.Item(97, 12).Value = 250000290  # L97
$ws.Cells.Item(97, 13).Value = -2194  # M97
$ws.Cells.Item(97, 14).Value = -250001282  # N97
$ws.Cells.Item(110, 8).Value = 721  # H110
$ws.Cells.Item(110, 9).Value = 721  # I110
$ws.Cells.Item(110, 10).Value = 0  # J110
$ws.Cells.Item(110, 11).Value = 721  # K110
$ws.Cells.Item(110, 12).Value = 0  # L110
$ws.Cells.Item(110, 13).Value = 1324  # M110
$ws.Cells.Item(110, 14).ClearContents()  # N110
$ws.Cells.Item(116, 8).Value = 1476.7222  # H116
$ws.Cells.Item(116, 9).Value = 1314.4375  # I116
$ws.Cells.Item(116, 11).Value = 1314.4375  # K116
$ws.Cells.Item(116, 13).Value = 979.5625  # M116
$ws.Cells.Item(132, 8).Value = 26938.143  # H132
$ws.Cells.Item(132, 9).Value = 3019.077  # I132
$ws.Cells.Item(132, 10).Value = 65806.625  # J132
$ws.Cells.Item(132, 11).Value = 9057.231  # K132
$ws.Cells.Item(132, 12).Value = 197419.875  # L132
$ws.Cells.Item(132, 13).Value = -6527.231  # M132
$ws.Cells.Item(132, 14).Value = -202479.875  # N132
$ws.Cells.Item(136, 8).Value = 3348.963  # H136
$ws.Cells.Item(136, 9).Value = 3216.88  # I136
$ws.Cells.Item(136, 10).Value = 5000  # J136
$ws.Cells.Item(136, 11).Value = 9650.639999999999  # K136
$ws.Cells.Item(136, 12).Value = 15000  # L136
$ws.Cells.Item(136, 13).Value = -7100.639999999999  # M136
$ws.Cells.Item(136, 14).Value = -20100  # N136

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1476.7222  # H3
$ws.Cells.Item(3, 9).Value = 1314.4375  # I3
$ws.Cells.Item(3, 11).Value = 1314.4375  # K3
$ws.Cells.Item(3, 13).Value = -1200.4375  # M3
$ws.Cells.Item(33, 8).Value = 6500  # H33
$ws.Cells.Item(33, 9).Value = 3000  # I33
$ws.Cells.Item(33, 10).Value = 10000  # J33
$ws.Cells.Item(33, 11).Value = 3000  # K33
$ws.Cells.Item(33, 12).Value = 10000  # L33
$ws.Cells.Item(33, 13).Value = -2664  # M33
$ws.Cells.Item(33, 14).Value = -10672  # N33
$ws.Cells.Item(94, 8).Value = 948.3889  # H94
$ws.Cells.Item(94, 9).Value = 771.75  # I94
$ws.Cells.Item(94, 11).Value = 771.75  # K94
$ws.Cells.Item(94, 13).Value = -320.75  # M94

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3130.3438  # H31
$ws.Cells.Item(31, 10).Value = 3525.577  # J31
$ws.Cells.Item(31, 12).Value = 3525.577  # L31
$ws.Cells.Item(31, 14).Value = -4115.577  # N31
$ws.Cells.Item(34, 8).Value = 3130.3438  # H34
$ws.Cells.Item(34, 10).Value = 3525.577  # J34
$ws.Cells.Item(34, 12).Value = 3525.577  # L34
$ws.Cells.Item(34, 14).Value = -3929.577  # N34
$ws.Cells.Item(38, 8).Value = 0  # H38
$ws.Cells.Item(38, 10).Value = 0  # J38
$ws.Cells.Item(38, 12).Value = 0  # L38
$ws.Cells.Item(38, 14).ClearContents()  # N38
$ws.Cells.Item(46, 8).Value = 0  # H46
$ws.Cells.Item(46, 10).Value = 0  # J46
$ws.Cells.Item(46, 12).Value = 0  # L46
$ws.Cells.Item(46, 14).ClearContents()  # N46
$ws.Cells.Item(59, 8).Value = 25880  # H59
$ws.Cells.Item(59, 10).Value = 25880  # J59
$ws.Cells.Item(59, 12).Value = 25880  # L59
$ws.Cells.Item(59, 14).Value = -28170  # N59
$ws.Cells.Item(127, 8).Value = 35113  # H127
$ws.Cells.Item(127, 10).Value = 35113  # J127
$ws.Cells.Item(127, 12).Value = 35113  # L127
$ws.Cells.Item(127, 14).Value = -45033  # N127
$ws.Cells.Item(134, 8).Value = 1240.25  # H134
$ws.Cells.Item(134, 9).Value = 1099.4375  # I134
$ws.Cells.Item(134, 11).Value = 3298.3125  # K134
$ws.Cells.Item(134, 13).Value = -763.3125  # M134

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(63, 8).Value = 4842.25  # H63
$ws.Cells.Item(63, 9).Value = 3549.75  # I63
$ws.Cells.Item(63, 10).Value = 6134.75  # J63
$ws.Cells.Item(63, 11).Value = 10649.25  # K63
$ws.Cells.Item(63, 12).Value = 18404.25  # L63
$ws.Cells.Item(63, 13).Value = -9900.25  # M63
$ws.Cells.Item(63, 14).Value = -19902.25  # N63
$ws.Cells.Item(66, 8).Value = 4842.25  # H66
$ws.Cells.Item(66, 9).Value = 3549.75  # I66
$ws.Cells.Item(66, 10).Value = 6134.75  # J66
$ws.Cells.Item(66, 11).Value = 31947.75  # K66
$ws.Cells.Item(66, 12).Value = 55212.75  # L66
$ws.Cells.Item(66, 13).Value = -28203.75  # M66
$ws.Cells.Item(66, 14).Value = -62700.75  # N66
$ws.Cells.Item(107, 8).Value = 3408.1614  # H107
$ws.Cells.Item(107, 9).Value = 6047.0586  # I107
$ws.Cells.Item(107, 10).Value = 203.78572  # J107
$ws.Cells.Item(107, 11).Value = 18141.1758  # K107
$ws.Cells.Item(107, 12).Value = 611.35716  # L107
$ws.Cells.Item(107, 13).Value = -16221.1758  # M107
$ws.Cells.Item(107, 14).Value = -4451.35716  # N107
$ws.Cells.Item(117, 8).Value = 1492.2  # H117
$ws.Cells.Item(117, 10).Value = 2032  # J117
$ws.Cells.Item(117, 12).Value = 6096  # L117
$ws.Cells.Item(117, 14).Value = -12980  # N117
$ws.Cells.Item(122, 8).Value = 480.7  # H122
$ws.Cells.Item(122, 9).Value = 368.83334  # I122
$ws.Cells.Item(122, 10).Value = 528.6429000000001  # J122
$ws.Cells.Item(122, 11).Value = 3319.50006  # K122
$ws.Cells.Item(122, 12).Value = 4757.7861  # L122
$ws.Cells.Item(122, 13).Value = -869.5000600000003  # M122
$ws.Cells.Item(122, 14).Value = -9657.786100000001  # N122
$ws.Cells.Item(131, 8).Value = 700.40405  # H131
$ws.Cells.Item(131, 10).Value = 719.23914  # J131
$ws.Cells.Item(131, 12).Value = 2157.71742  # L131
$ws.Cells.Item(131, 14).Value = -12237.71742  # N131
$ws.Cells.Item(140, 8).Value = 2837.9546  # H140
$ws.Cells.Item(140, 9).Value = 1486  # I140
$ws.Cells.Item(140, 10).Value = 3964.5833  # J140
$ws.Cells.Item(140, 11).Value = 4458  # K140
$ws.Cells.Item(140, 12).Value = 11893.7499  # L140
$ws.Cells.Item(140, 13).Value = 722  # M140
$ws.Cells.Item(140, 14).Value = -22253.7499  # N140

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 69.36364  # H2
$ws.Cells.Item(2, 9).Value = 62.875  # I2
$ws.Cells.Item(2, 10).Value = 86.666664  # J2
$ws.Cells.Item(2, 11).Value = 62.875  # K2
$ws.Cells.Item(2, 12).Value = 86.666664  # L2
$ws.Cells.Item(2, 13).Value = 50.125  # M2
$ws.Cells.Item(2, 14).Value = -312.666664  # N2

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 702.9091  # H16
$ws.Cells.Item(16, 9).Value = 653.2  # I16
$ws.Cells.Item(16, 11).Value = 653.2  # K16
$ws.Cells.Item(16, 13).Value = -483.2  # M16
$ws.Cells.Item(36, 8).Value = 35000  # H36
$ws.Cells.Item(36, 10).Value = 35000  # J36
$ws.Cells.Item(36, 12).Value = 35000  # L36
$ws.Cells.Item(36, 14).Value = -36124  # N36
$ws.Cells.Item(100, 8).Value = 2250.75  # H100
$ws.Cells.Item(100, 9).Value = 1854  # I100
$ws.Cells.Item(100, 10).Value = 2383  # J100
$ws.Cells.Item(100, 11).Value = 1854  # K100
$ws.Cells.Item(100, 12).Value = 2383  # L100
$ws.Cells.Item(100, 13).Value = -1313  # M100
$ws.Cells.Item(100, 14).Value = -3465  # N100

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 15214.833  # H2
$ws.Cells.Item(2, 10).Value = 20500.25  # J2
$ws.Cells.Item(2, 12).Value = 20500.25  # L2
$ws.Cells.Item(2, 14).Value = -20724.25  # N2
$ws.Cells.Item(100, 8).Value = 319.8  # H100
$ws.Cells.Item(100, 9).Value = 333.1111  # I100
$ws.Cells.Item(100, 11).Value = 666.2222  # K100
$ws.Cells.Item(100, 13).Value = -125.2222  # M100
$ws.Cells.Item(113, 8).Value = 2159.1538  # H113
$ws.Cells.Item(113, 9).Value = 2508.7273  # I113
$ws.Cells.Item(113, 10).Value = 236.5  # J113
$ws.Cells.Item(113, 11).Value = 7526.1819  # K113
$ws.Cells.Item(113, 12).Value = 709.5  # L113
$ws.Cells.Item(113, 13).Value = -5356.1819  # M113
$ws.Cells.Item(113, 14).Value = -5049.5  # N113
$ws.Cells.Item(132, 8).Value = 1389.8572  # H132
$ws.Cells.Item(132, 9).Value = 800.8570999999999  # I132
$ws.Cells.Item(132, 11).Value = 2402.5713  # K132
$ws.Cells.Item(132, 13).Value = 127.4287000000004  # M132
$ws.Cells.Item(136, 8).Value = 21068856  # H136
$ws.Cells.Item(136, 9).Value = 27165684  # I136
$ws.Cells.Item(136, 10).Value = 7091.364  # J136
$ws.Cells.Item(136, 11).Value = 81497052  # K136
$ws.Cells.Item(136, 12).Value = 21274.092  # L136
$ws.Cells.Item(136, 13).Value = -81494502  # M136
$ws.Cells.Item(136, 14).Value = -26374.092  # N136
